$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "41.526.22"
$ws.Range("E2").Value = "  +0.71%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.473.46"
$ws.Range("E3").Value = "  +0.48%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  +0.11%  "

# Row 6 - Solana
$ws.Range("D6").Value = "91.63"
$ws.Range("E6").Value = "  -2.48%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.23%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  +2.70%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "32.48"
$ws.Range("E10").Value = "  -2.63%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  +0.97%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +1.46%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.856.60"
$ws.Range("E13").Value = "  +0.57%  "

# Row 14 - was Chainlink, now Polkadot
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "6.86"
$ws.Range("E14").Value = "  -1.52%  "

# Row 15 - was Polkadot, now Chainlink
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "16.25"
$ws.Range("E15").Value = "  +8.98%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.464.93"
$ws.Range("E16").Value = "  -0.25%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.770"
$ws.Range("E17").Value = "  -1.92%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "41.517.23"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "6.51"
$ws.Range("E19").Value = "  +3.29%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  +1.73%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "71.81"
$ws.Range("E21").Value = "  +4.97%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("D22").Value = "11.02"
$ws.Range("E22").Value = "  -1.45%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "236.15"
$ws.Range("E23").Value = "  -0.42%  "

# Row 24 - PancakeSwap
$ws.Range("D24").Value = "2.71"
$ws.Range("E24").Value = "  -1.35%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  -0.11%  "

# Row 26 - ImmutableX
$ws.Range("D26").Value = "1.89"
$ws.Range("E26").Value = "  +0.23%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "24.93"
$ws.Range("E27").Value = "  +4.27%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  -0.56%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "9.67"
$ws.Range("E29").Value = "  +0.29%  "

# Row 30 - InjectiveProtocol
$ws.Range("D30").Value = "35.65"
$ws.Range("E30").Value = "  -1.85%  "

# Row 31 - Monero
$ws.Range("D31").Value = "157.12"
$ws.Range("E31").Value = "  +3.70%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "5.43"
$ws.Range("E32").Value = "  -0.73%  "

# Row 33 - WEMIXToken
$ws.Range("E33").Value = "  -0.53%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.0756"
$ws.Range("E34").Value = "  +1.63%  "

# Row 35 - Celestia
$ws.Range("D35").Value = "17.38"
$ws.Range("E35").Value = "  +1.52%  "

# Row 36 - ApeXProtocol
$ws.Range("D36").Value = "2.39"
$ws.Range("E36").Value = "  -8.33%  "

# Row 37 - was Kaspa, now LidoDAOToken
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "2.89"
$ws.Range("E37").Value = "  -5.15%  "

# Row 38 - was LidoDAOToken, now Kaspa
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.104"
$ws.Range("E38").Value = "  +2.88%  "

# Row 39 - ARBITRUM
$ws.Range("D39").Value = "1.81"
$ws.Range("E39").Value = "  -2.59%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  -0.12%  "

# Row 41 - RenderToken
$ws.Range("E41").Value = "  -4.72%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  -0.24%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.959.21"
$ws.Range("E43").Value = "  -1.08%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -0.27%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "18.89"
$ws.Range("E45").Value = "  -2.80%  "

# Row 46 - NEARProtocol
$ws.Range("E46").Value = "  -2.57%  "

# Row 47 - FraxShare
$ws.Range("D47").Value = "8.97"
$ws.Range("E47").Value = "  +3.34%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "2.717.25"

# Row 49 - Aave
$ws.Range("D49").Value = "97.54"
$ws.Range("E49").Value = "  +1.30%  "

# Row 50 - ordi
$ws.Range("D50").Value = "67.32"
$ws.Range("E50").Value = "  -3.03%  "

# Row 51 - BitcoinSV
$ws.Range("D51").Value = "72.04"
$ws.Range("E51").Value = "  -3.18%  "
